# Apply the "SVM V3 with 82% accuracy" edits to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update overall accuracy value in B1 (90 -> 82)
$ws.Range("B1").Value = 82

# Update individual prediction/actual cells (B = predicted, C = actual columns)
$ws.Range("C5").Value = 0    # BSD 03

$ws.Range("B7").Value = 0    # BSD 04
$ws.Range("C7").Value = 0    # BSD 04

$ws.Range("C10").Value = 1   # BSD 07

$ws.Range("C12").Value = 1   # BSD 09

$ws.Range("B13").Value = 1   # BSD 10

$ws.Range("C16").Value = 0   # BSD 13

$ws.Range("B20").Value = 0   # BSD 17

$ws.Range("B26").Value = 0   # BSD 23
$ws.Range("C26").Value = 0   # BSD 23

$ws.Range("B27").Value = 0   # BSD 24

$ws.Range("C31").Value = 0   # BSD 28

$ws.Range("C32").Value = 1   # BSD 29

# Update the active selection to match the saved view state (F28)
$ws.Range("F28").Select()
